# The author revised two quiz-question prompts (adding clarifying
# sentences) in the "fins.xlsx" workbook, then left the workbook with a
# few different cell selections / a different active sheet than before.
#
# Sheet tab names -> internal file mapping (for reference):
#   "0_" "1_" "2_" "3_" "4_" "5_" "6_" "Sheet1"

$wb = $excel.ActiveWorkbook

# --- Edit 1: clarify the "large P/A ratio" question on sheet "1_" (cell A1) ---
$shPA = $wb.Worksheets.Item("1_")
$shPA.Activate()
$shPA.Range("A1").Value = 'A large second derivative implies a "cooler" fin tip.  A large  "P/A" ratio increases the size of the derivative in the fin equation.  Why is this true, in physical terms?  Remember that "A" is the cross-sectional area and not the surface area (Mark all that are true)'
$shPA.Range("B2:B3").Select()

# --- Edit 2: clarify the "adiabatic tip" question on sheet "3_" (cell A1) ---
$shAdiabatic = $wb.Worksheets.Item("3_")
$shAdiabatic.Activate()
$shAdiabatic.Range("A1").Value = 'An adiabatic (no flux) BC at the tip suggests that flux through the tip *in the x-direction* is insignificant: that is, all of the thermal energy that enters the fin leaves through the sides.  Think about what the temperature field would look like in this case (what would dT/dx be?), and determine the "flux_base" that leads to this condition.  Answer as a multiple of 10^5 W/m2, accurate to 2 decimal places'
$shAdiabatic.Range("A5:A6").Select()

# --- The user also browsed sheet "6_" and left a cell selected there ---
$shLast = $wb.Worksheets.Item("6_")
$shLast.Activate()
$shLast.Range("B3").Select()

# --- Finish back on sheet "3_", which is the sheet left active on save ---
$shAdiabatic.Activate()
